# ========================================================================
# results.xlsx edit: "update predict.py and color npm dataset"
#   - rename sheet3 semantic2npm -> common_class
#   - sheet1 (semantic3d): drop stray row 16, finish row 9, add row 10
#   - sheet2 (npm3d): add row 3
#   - sheet3 (common_class): insert a "dataset_transfer" column, add row 3
# ========================================================================

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

$ws3.Name = "common_class"

# ------------------------------------------------------------------
# Sheet1 (semantic3d)
# ------------------------------------------------------------------
$ws1.Rows.Item(16).Delete() | Out-Null

$ws1.Range("I9").Value = 40
$ws1.Range("J9").Value = 0.65811299999999995
$ws1.Range("K9").Value = 0.89662500000000001
$ws1.Range("L9").Value = 0.90305599999999997
$ws1.Range("M9").Value = 0.69113800000000003
$ws1.Range("N9").Value = 0.895814
$ws1.Range("O9").Value = 0.37152800000000002
$ws1.Range("P9").Value = 0.91305499999999995
$ws1.Range("Q9").Value = 0.316799
$ws1.Range("R9").Value = 0.55424600000000002
$ws1.Range("S9").Value = 0.61926800000000004

$ws1.Range("A10").Value = "semantic_cross_npm"
$ws1.Range("B10").Value = 8192
$ws1.Range("C10").Value = 16
$ws1.Range("D10").Value = 1
$ws1.Range("E10").Value = 0
$ws1.Range("F10").Value = 0
$ws1.Range("G10").Value = 10
$ws1.Range("H10").Value = 10
$ws1.Range("I10").Value = 1
$ws1.Range("J10").Value = 0.13391500000000001
$ws1.Range("K10").Value = 0.44855600000000001
$ws1.Range("L10").Value = 0.38561000000000001
$ws1.Range("M10").Value = 0.237426
$ws1.Range("N10").Value = 0.000286
$ws1.Range("O10").Value = 0
$ws1.Range("P10").Value = 0.43162800000000001
$ws1.Range("Q10").Value = 0.016371
$ws1.Range("R10").Value = 0
$ws1.Range("S10").Value = 0

$ws1.Columns.Item(1).ColumnWidth = 18.910714285714285
$ws1.Columns.Item(12).ColumnWidth = 11.910714285714286

$ws1.Range("G23").Select() | Out-Null

# ------------------------------------------------------------------
# Sheet2 (npm3d)
# ------------------------------------------------------------------
$ws2.Range("A3").Value = "ours"
$ws2.Range("B3").Value = 8192
$ws2.Range("C3").Value = 16
$ws2.Range("D3").Value = 1
$ws2.Range("E3").Value = 0
$ws2.Range("F3").Value = 0
$ws2.Range("G3").Value = 10
$ws2.Range("H3").Value = 10
$ws2.Range("I3").Value = 50
$ws2.Range("J3").Value = 0.66985799999999995
$ws2.Range("K3").Value = 0.962615
$ws2.Range("L3").Value = 0.97357199999999999
$ws2.Range("M3").Value = 0.95247999999999999
$ws2.Range("N3").Value = 0.64186900000000002
$ws2.Range("O3").Value = 0.59559300000000004
$ws2.Range("P3").Value = 0.40692699999999998
$ws2.Range("Q3").Value = 0.40978100000000001
$ws2.Range("R3").Value = 0.39611200000000002
$ws2.Range("S3").Value = 0.84192199999999995
$ws2.Range("T3").Value = 0.81046499999999999

$ws2.Range("E17").Select() | Out-Null

# ------------------------------------------------------------------
# Sheet3 (common_class) - insert new "dataset_transfer" column at B
# ------------------------------------------------------------------
$ws3.Columns.Item(2).Insert() | Out-Null

$ws3.Range("B1").Value = "dataset_transfer"
$ws3.Range("B2").Value = "semantic2npm"

$ws3.Range("A3").Value = "pointsemantic"
$ws3.Range("B3").Value = "npm2npm"
$ws3.Range("C3").Value = 2
$ws3.Range("D3").Value = 0.83209999999999995
$ws3.Range("E3").Value = 0.96519999999999995
$ws3.Range("F3").Value = 0.97370000000000001
$ws3.Range("G3").Value = 0.81569999999999998
$ws3.Range("H3").Value = 0.95469999999999999
$ws3.Range("I3").Value = 0.52929999999999999
$ws3.Range("J3").Value = 0.88680000000000003

$ws3.Columns.Item(2).ColumnWidth = 14.535714285714286

$ws3.Range("F11").Select() | Out-Null

Write-Output "edits applied"
